$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of EUR->ARS quote data. Column A holds a date-looking string
# ("2025-09-15") that Excel would otherwise auto-convert to a date serial,
# so force it to text via NumberFormat, assign, then restore the cell
# style so no stray formatting is left behind (matches sibling rows).
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2025-09-15"
$ws.Range("A20").Style = "Normal"

$ws.Range("B20").Value = "21:20:36"
$ws.Range("C20").Value = "1.00 EUR = 1,712.3015"
